# Weekly update: insert a new price observation as row 22, shifting the
# existing rows 22-82 down to 23-83 (dimension grows from R82 to R83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 22; Excel shifts rows 22..82 down to 23..83
# and copies the formatting (incl. the date style on column D) from the row
# that was previously at 22.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly observation.
$ws.Cells.Item(22, 1).Value  = 9
$ws.Cells.Item(22, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(22, 3).Value  = "Metropolitana"
$ws.Cells.Item(22, 4).Value  = 44526
$ws.Cells.Item(22, 5).Value  = 13
$ws.Cells.Item(22, 6).Value  = 100112022
$ws.Cells.Item(22, 7).Value  = "Arveja Verde"
$ws.Cells.Item(22, 8).Value  = "Sin especificar"
$ws.Cells.Item(22, 9).Value  = "Primera"
$ws.Cells.Item(22, 10).Value = 25
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 16000
$ws.Cells.Item(22, 13).Value = 15480
$ws.Cells.Item(22, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(22, 15).Value = "Carahue"
$ws.Cells.Item(22, 16).Value = 619
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"
